# Apply the "Added periodic & upfront related scenarios" change.
#
# Semantically this edit:
#   1. Changes the "repaymentstrategy" value on the ProductLoanInput sheet
#      (cell B17) from "RBI (India)" to "Overdue/Due Fee/Int,Principal".
#   2. Updates the sheet view so that the selected / visible cell follows
#      the edited cell (topLeftCell A7, active cell B17) instead of the
#      previous selection near the bottom of the sheet.
#
# (The removal of the now-unused "RBI (India)" shared string and the
# associated shift of every other shared-string index is just a side
# effect of Excel rebuilding the shared strings table when the workbook
# is saved - it carries no additional semantic change.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")

$ws.Range("B17").Value = "Overdue/Due Fee/Int,Principal"

$ws.Activate()
$ws.Range("B17").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
